$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at position 38 (AL), pushing existing AL:AR right to AM:AS
$ws.Columns.Item(38).Insert()

# New column AL (38) header text in row 10 and value in row 11
$ws.Cells.Item(10, 38).Value = "{bgPensumZeiteinheitTitle}"
$ws.Cells.Item(11, 38).Value = "{bgPensumZeiteinheit}"
